# Apply PCSMOTE log corrections: percentil_densidad_25 (G) threshold fix,
# recomputed percentil_dist_75 (F) and densidad (K) values, refreshed timestamps (Z).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 112

# New percentil_dist_75 values (column F), one entry per data row (row 2..112)
$fVals = @(1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.003449896865958,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,1.000000011720752,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9886747427800949,0.9465529195601421,0.9465529195601421,0.9465529195601421,0.9465529195601421,0.9465529195601421,0.9465529195601421,0.9465529195601421,0.9465529195601421,0.9465529195601421,0.9465529195601421)

# New densidad values (column K), one entry per data row (row 2..112)
$kVals = @(0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,1,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.4285714285714285,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.7142857142857143,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857,0.2857142857142857)

# Refreshed timestamp values (column Z), one entry per data row (row 2..112)
$zVals = @("2025-10-19T23:56:01.994238","2025-10-19T23:56:01.995212","2025-10-19T23:56:01.995212","2025-10-19T23:56:01.995212","2025-10-19T23:56:01.995212","2025-10-19T23:56:01.995212","2025-10-19T23:56:01.996212","2025-10-19T23:56:01.996212","2025-10-19T23:56:01.996212","2025-10-19T23:56:01.996212","2025-10-19T23:56:01.996212","2025-10-19T23:56:01.997210","2025-10-19T23:56:01.998208","2025-10-19T23:56:01.998208","2025-10-19T23:56:01.999213","2025-10-19T23:56:01.999213","2025-10-19T23:56:01.999213","2025-10-19T23:56:01.999213","2025-10-19T23:56:01.999213","2025-10-19T23:56:02.000214","2025-10-19T23:56:02.000214","2025-10-19T23:56:02.000214","2025-10-19T23:56:02.000214","2025-10-19T23:56:02.001212","2025-10-19T23:56:02.001212","2025-10-19T23:56:02.001212","2025-10-19T23:56:02.001212","2025-10-19T23:56:02.001212","2025-10-19T23:56:02.002213","2025-10-19T23:56:02.002213","2025-10-19T23:56:02.002213","2025-10-19T23:56:02.002213","2025-10-19T23:56:02.002213","2025-10-19T23:56:02.003214","2025-10-19T23:56:02.003214","2025-10-19T23:56:02.003214","2025-10-19T23:56:02.003214","2025-10-19T23:56:02.004211","2025-10-19T23:56:02.004211","2025-10-19T23:56:02.004211","2025-10-19T23:56:02.004211","2025-10-19T23:56:02.004211","2025-10-19T23:56:02.005213","2025-10-19T23:56:02.005213","2025-10-19T23:56:02.036208","2025-10-19T23:56:02.036208","2025-10-19T23:56:02.036208","2025-10-19T23:56:02.037210","2025-10-19T23:56:02.037210","2025-10-19T23:56:02.037210","2025-10-19T23:56:02.038208","2025-10-19T23:56:02.038208","2025-10-19T23:56:02.038208","2025-10-19T23:56:02.038208","2025-10-19T23:56:02.039212","2025-10-19T23:56:02.039212","2025-10-19T23:56:02.039212","2025-10-19T23:56:02.039212","2025-10-19T23:56:02.039212","2025-10-19T23:56:02.040209","2025-10-19T23:56:02.040209","2025-10-19T23:56:02.040209","2025-10-19T23:56:02.040209","2025-10-19T23:56:02.041312","2025-10-19T23:56:02.042396","2025-10-19T23:56:02.042396","2025-10-19T23:56:02.042396","2025-10-19T23:56:02.043380","2025-10-19T23:56:02.043380","2025-10-19T23:56:02.043380","2025-10-19T23:56:02.043380","2025-10-19T23:56:02.043380","2025-10-19T23:56:02.044375","2025-10-19T23:56:02.077107","2025-10-19T23:56:02.078107","2025-10-19T23:56:02.078107","2025-10-19T23:56:02.078107","2025-10-19T23:56:02.078107","2025-10-19T23:56:02.079107","2025-10-19T23:56:02.079107","2025-10-19T23:56:02.079107","2025-10-19T23:56:02.079107","2025-10-19T23:56:02.079107","2025-10-19T23:56:02.080107","2025-10-19T23:56:02.080107","2025-10-19T23:56:02.081105","2025-10-19T23:56:02.081105","2025-10-19T23:56:02.082107","2025-10-19T23:56:02.082107","2025-10-19T23:56:02.082107","2025-10-19T23:56:02.083109","2025-10-19T23:56:02.083109","2025-10-19T23:56:02.083109","2025-10-19T23:56:02.083109","2025-10-19T23:56:02.083109","2025-10-19T23:56:02.084106","2025-10-19T23:56:02.084106","2025-10-19T23:56:02.084106","2025-10-19T23:56:02.084106","2025-10-19T23:56:02.084106","2025-10-19T23:56:02.084106","2025-10-19T23:56:02.111104","2025-10-19T23:56:02.111104","2025-10-19T23:56:02.111104","2025-10-19T23:56:02.111104","2025-10-19T23:56:02.111104","2025-10-19T23:56:02.112104","2025-10-19T23:56:02.112104","2025-10-19T23:56:02.112104","2025-10-19T23:56:02.112104","2025-10-19T23:56:02.113103")

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $idx = $row - $firstRow
    $ws.Cells.Item($row, 6).Value = $fVals[$idx]
    $ws.Cells.Item($row, 7).Value = 0.2857142857142857
    $ws.Cells.Item($row, 11).Value = $kVals[$idx]
    $ws.Cells.Item($row, 26).Value = $zVals[$idx]
}

Write-Host "Updated percentil_dist_75, percentil_densidad_25, densidad and timestamp for rows $firstRow..$lastRow"
